# Rename the "Plain_English" column header (underscore) to "Plain English"
# (space) across the classification/"to_Code" lookup sheets, matching the
# latest RD model naming convention.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Fuel_to_Code").Range("C1").Value = "Plain English"
$wb.Worksheets.Item("VehFuel_to_Code").Range("C1").Value = "Plain English"
$wb.Worksheets.Item("Tech_to_Code").Range("C1").Value = "Plain English"
$wb.Worksheets.Item("Dem_to_Code").Range("B1").Value = "Plain English"
